# The codeforIATI SectorGroup codelist re-ordered the four
# "codeforiati:*" columns: column D swaps places with column F, and
# column E swaps places with column G. (E.g. header D used to read
# "codeforiati:category-code" and now reads "codeforiati:group-name";
# a data row's D used to hold the category-code value "112" and now
# holds the group-name value "Education", while F/G pick up what D/E
# used to hold.) This applies uniformly to every row, including the
# header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVal = [string]$ws.Cells.Item($r, 4).Text
    $eVal = [string]$ws.Cells.Item($r, 5).Text
    $fVal = [string]$ws.Cells.Item($r, 6).Text
    $gVal = [string]$ws.Cells.Item($r, 7).Text

    # Every one of these columns is stored as text (codes like "111"/"110"
    # as well as plain names), so force text format before writing back to
    # avoid Excel auto-converting numeric-looking strings to real numbers.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 7).NumberFormat = "@"

    # Column D swaps with column F, column E swaps with column G
    $ws.Cells.Item($r, 4).Value2 = $fVal
    $ws.Cells.Item($r, 5).Value2 = $gVal
    $ws.Cells.Item($r, 6).Value2 = $dVal
    $ws.Cells.Item($r, 7).Value2 = $eVal
}
